$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 120
$ws.Range("A120").Value = 4
$ws.Range("B120").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C120").Value = "Los Lagos"
$ws.Range("D120").Value = 44939
$ws.Range("D120").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E120").Value = 10
$ws.Range("F120").Value = "Fruta"
$ws.Range("G120").Value = 100103
$ws.Range("H120").Value = "Frutos de hueso (carozo)"
$ws.Range("I120").Value = 100103001
$ws.Range("J120").Value = "Cereza"
$ws.Range("K120").Value = "Bing"
$ws.Range("L120").Value = "Primera"
$ws.Range("M120").Value = 800
$ws.Range("N120").Value = 5000
$ws.Range("O120").Value = 5500
$ws.Range("P120").Value = 5250
$ws.Range("Q120").Value = "`$/bandeja 5 kilos"
$ws.Range("R120").Value = "Provincia de Curicó"
$ws.Range("S120").Value = 1050
$ws.Range("T120").Value = 5

# New row 121
$ws.Range("A121").Value = 4
$ws.Range("B121").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C121").Value = "Los Lagos"
$ws.Range("D121").Value = 44939
$ws.Range("D121").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E121").Value = 10
$ws.Range("F121").Value = "Fruta"
$ws.Range("G121").Value = 100103
$ws.Range("H121").Value = "Frutos de hueso (carozo)"
$ws.Range("I121").Value = 100103001
$ws.Range("J121").Value = "Cereza"
$ws.Range("K121").Value = "Lapins"
$ws.Range("L121").Value = "Primera"
$ws.Range("M121").Value = 600
$ws.Range("N121").Value = 5000
$ws.Range("O121").Value = 5500
$ws.Range("P121").Value = 5250
$ws.Range("Q121").Value = "`$/bandeja 5 kilos"
$ws.Range("R121").Value = "Provincia de Curicó"
$ws.Range("S121").Value = 1050
$ws.Range("T121").Value = 5
